$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Write the combined header + canton data block (A1:K24) ---
$arr = New-Object 'object[,]' 24,11
$arr[0,0] = "idx"
$arr[0,1] = "idx2"
$arr[0,2] = "Name"
$arr[0,3] = "Date Start"
$arr[0,4] = "Date End"
$arr[0,5] = "(m3/s)"
$arr[0,6] = "(MW1)"
$arr[0,7] = "(MW2)"
$arr[0,8] = "(GWh) Winter"
$arr[0,9] = "(GWh) Summer"
$arr[0,10] = "(GWh) Year"
$arr[1,0] = 1
$arr[1,1] = 509300
$arr[1,2] = "Taulan"
$arr[1,3] = 1887
$arr[1,4] = 1996
$arr[1,5] = 0.68
$arr[1,6] = 1.2
$arr[1,7] = 1.2
$arr[1,8] = 1.5
$arr[1,9] = 1.9
$arr[1,10] = 3.4
$arr[2,0] = 2
$arr[2,1] = 205600
$arr[2,2] = "Le Chalet"
$arr[2,3] = 1894
$arr[2,4] = 1988
$arr[2,5] = 18
$arr[2,6] = 2.6
$arr[2,7] = 2.3199999999999998
$arr[2,8] = 5.7
$arr[2,9] = 3.6
$arr[2,10] = 9.3000000000000007
$arr[3,0] = 3
$arr[3,1] = 203800
$arr[3,2] = "Montbovon"
$arr[3,3] = 1896
$arr[3,4] = 1972
$arr[3,5] = 40
$arr[3,6] = 24.08
$arr[3,7] = 22.14
$arr[3,8] = 13.98
$arr[3,9] = 47.38
$arr[3,10] = 61.37
$arr[4,0] = 4
$arr[4,1] = 509400
$arr[4,2] = "Plan-Dessous"
$arr[4,3] = 1896
$arr[4,4] = 2001
$arr[4,5] = 10.7
$arr[4,6] = 10.7
$arr[4,7] = 7.85
$arr[4,8] = 16.8
$arr[4,9] = 10.7
$arr[4,10] = 27.5
$arr[5,0] = 5
$arr[5,1] = 508100
$arr[5,2] = "Sublin 1"
$arr[5,3] = 1898
$arr[5,4] = 1993
$arr[5,5] = 5.2
$arr[5,6] = 8
$arr[5,7] = 7.2
$arr[5,8] = 8
$arr[5,9] = 23
$arr[5,10] = 31
$arr[6,0] = 6
$arr[6,1] = 509200
$arr[6,2] = "Sonzier"
$arr[6,3] = 1901
$arr[6,4] = 1971
$arr[6,5] = 0.5
$arr[6,6] = 1.65
$arr[6,7] = 1.6
$arr[6,8] = 2
$arr[6,9] = 4.5999999999999996
$arr[6,10] = 6.6
$arr[7,0] = 7
$arr[7,1] = 205200
$arr[7,2] = "La Dernier"
$arr[7,3] = 1903
$arr[7,4] = 1988
$arr[7,5] = 13
$arr[7,6] = 28
$arr[7,7] = 27
$arr[7,8] = 22
$arr[7,9] = 9
$arr[7,10] = 31
$arr[8,0] = 8
$arr[8,1] = 508900
$arr[8,2] = "Les Farettes"
$arr[8,3] = 1906
$arr[8,4] = 1967
$arr[8,5] = 2.5
$arr[8,6] = 15
$arr[8,7] = 6.7
$arr[8,8] = 26
$arr[8,9] = 31
$arr[8,10] = 57
$arr[9,0] = 9
$arr[9,1] = 205500
$arr[9,2] = "Montcherand"
$arr[9,3] = 1908
$arr[9,4] = 1950
$arr[9,5] = 19
$arr[9,6] = 15.5
$arr[9,7] = 14
$arr[9,8] = 26.5
$arr[9,9] = 30.5
$arr[9,10] = 57
$arr[10,0] = 10
$arr[10,1] = 508200
$arr[10,2] = "Sublin 2"
$arr[10,3] = 1911
$arr[10,4] = 2002
$arr[10,5] = 0.25
$arr[10,6] = 1.5
$arr[10,7] = 1.3
$arr[10,8] = 2
$arr[10,9] = 4.0999999999999996
$arr[10,10] = 6.1
$arr[11,0] = 11
$arr[11,1] = 508800
$arr[11,2] = "Pont de la Tine"
$arr[11,3] = 1913
$arr[11,4] = 1991
$arr[11,5] = 2.5
$arr[11,6] = 10.6
$arr[11,7] = 5.3
$arr[11,8] = 14
$arr[11,9] = 19.600000000000001
$arr[11,10] = 33.6
$arr[12,0] = 12
$arr[12,1] = 205700
$arr[12,2] = "Moulins ROD"
$arr[12,3] = 1920
$arr[12,4] = 1982
$arr[12,5] = 15
$arr[12,6] = 0.57999999999999996
$arr[12,7] = 0.57999999999999996
$arr[12,8] = 0.9
$arr[12,9] = 0.8
$arr[12,10] = 1.7
$arr[13,0] = 13
$arr[13,1] = 508000
$arr[13,2] = "La Peuffeyre"
$arr[13,3] = 1927
$arr[13,4] = 2004
$arr[13,5] = 6.6
$arr[13,6] = 24
$arr[13,7] = 22
$arr[13,8] = 26.6
$arr[13,9] = 45.4
$arr[13,10] = 72
$arr[14,0] = 14
$arr[14,1] = 508300
$arr[14,2] = "Bevieux"
$arr[14,3] = 1943
$arr[14,4] = $null
$arr[14,5] = 4.0999999999999996
$arr[14,6] = 1.9
$arr[14,7] = 1.8
$arr[14,8] = 4.2
$arr[14,9] = 7.3
$arr[14,10] = 11.5
$arr[15,0] = 15
$arr[15,1] = 507900
$arr[15,2] = "Lavey"
$arr[15,3] = 1950
$arr[15,4] = 1990
$arr[15,5] = 220
$arr[15,6] = 37.799999999999997
$arr[15,7] = 29.4
$arr[15,8] = 79.8
$arr[15,9] = 88.2
$arr[15,10] = 168
$arr[16,0] = 16
$arr[16,1] = 205400
$arr[16,2] = "Les Clées"
$arr[16,3] = 1955
$arr[16,4] = $null
$arr[16,5] = 21
$arr[16,6] = 30
$arr[16,7] = 27
$arr[16,8] = 47.5
$arr[16,9] = 55.5
$arr[16,10] = 103
$arr[17,0] = 17
$arr[17,1] = 205300
$arr[17,2] = "La Jougnenaz"
$arr[17,3] = 1955
$arr[17,4] = 1970
$arr[17,5] = 6
$arr[17,6] = 2.35
$arr[17,7] = 2.1
$arr[17,8] = 3.8
$arr[17,9] = 2.2000000000000002
$arr[17,10] = 6
$arr[18,0] = 18
$arr[18,1] = 508700
$arr[18,2] = "Diablerets"
$arr[18,3] = 1957
$arr[18,4] = $null
$arr[18,5] = 1.75
$arr[18,6] = 2.0099999999999998
$arr[18,7] = 1.93
$arr[18,8] = 3.5
$arr[18,9] = 2.16
$arr[18,10] = 5.65
$arr[19,0] = 19
$arr[19,1] = 509100
$arr[19,2] = "Veytaux"
$arr[19,3] = 1972
$arr[19,4] = $null
$arr[19,5] = 32.6
$arr[19,6] = 146.06
$arr[19,7] = 146.06
$arr[19,8] = 60.25
$arr[19,9] = 52.95
$arr[19,10] = 113.2
$arr[20,0] = 20
$arr[20,1] = 508850
$arr[20,2] = "Douve 1"
$arr[20,3] = 1989
$arr[20,4] = 2000
$arr[20,5] = 0.1
$arr[20,6] = 0.46
$arr[20,7] = 0.42
$arr[20,8] = 1.1000000000000001
$arr[20,9] = 1
$arr[20,10] = 2.1
$arr[21,0] = 21
$arr[21,1] = 203750
$arr[21,2] = "Gérignoz (La Pontia)"
$arr[21,3] = 1996
$arr[21,4] = $null
$arr[21,5] = 0.19
$arr[21,6] = 0.45
$arr[21,7] = 0.45
$arr[21,8] = 1.41
$arr[21,9] = 0.94
$arr[21,10] = 2.35
$arr[22,0] = 22
$arr[22,1] = 508950
$arr[22,2] = "Fontanney"
$arr[22,3] = 1997
$arr[22,4] = $null
$arr[22,5] = 0.41
$arr[22,6] = 0.33
$arr[22,7] = 0.31
$arr[22,8] = 0.3
$arr[22,9] = 0.5
$arr[22,10] = 0.8
$arr[23,0] = 23
$arr[23,1] = 509425
$arr[23,2] = "La Petite Vaux"
$arr[23,3] = 2008
$arr[23,4] = $null
$arr[23,5] = 10
$arr[23,6] = 3.31
$arr[23,7] = 3.19
$arr[23,8] = 6.67
$arr[23,9] = 4.4400000000000004
$arr[23,10] = 11.11
$ws.Range("A1:K24").Value = $arr

# --- Formatting: new header row 1 (A1:E1 stay default Normal style; F1:K1 match the small-font label style) ---
$ws.Range("F1:K1").Font.Size = 9

# --- Formatting: data rows 2-24 ---
$ws.Range("A2:B24").Font.Size = 9
$ws.Range("A2:B24").NumberFormat = "0"
$ws.Range("D2:E24").Font.Size = 9
$ws.Range("D2:E24").NumberFormat = "0"
$ws.Range("C2:C24").Font.Size = 9
$ws.Range("F2:K24").Font.Size = 9
$ws.Range("F2:K24").NumberFormat = "0.00"

# --- New empty filler rows 25:26 (match style used by row 27 onward) ---
$ws.Range("A25:K26").ClearContents()
$ws.Range("A25:K26").Font.Size = 9
$ws.Range("A25:K26").NumberFormat = "0"

# --- Rows 83:84 no longer carry A:K filler cells ---
$ws.Range("A83:K84").Clear()

# --- Update the active selection to match the author edit ---
$ws.Range("A21:K21").Select()
